$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Expand the table to include the new column first (Resize auto-names the new
# column "Column4" as a placeholder header, so the real header text must be
# written afterwards or it gets clobbered).
$table = $ws.ListObjects.Item("States_and_Territories")
$table.Resize($ws.Range("A1:D10"))

# Add the new column header
$ws.Range("D1").Value = "Geography Sort Order"

# Fill the sort-order values for the 9 data rows
$sortOrders = 1..9
for ($i = 0; $i -lt $sortOrders.Count; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $sortOrders[$i]
}

# Column widths to match the authored (best-fit) layout as closely as possible
$ws.Columns.Item(1).ColumnWidth = 16.09
$ws.Columns.Item(3).ColumnWidth = 22.25
$ws.Columns.Item(4).ColumnWidth = 20.42

# Match the saved selection
$ws.Range("D9").Select()
